# "Add files via upload" — the re-uploaded worksheet has every date in
# column A (rows 2-115 of the "lista" sheet) advanced to the next entry
# in the list: row N's date becomes the date that used to be in row N+1,
# and the final row picks up a brand-new date (one day after the old
# last date). Columns B:K are untouched.
#
# Reproduce that by reading the current column-A values, then writing
# the shifted sequence back so each cell ends up holding what used to be
# its successor's value (with the very last cell getting "last + 1").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 115
$col = 1  # column A

# Snapshot the existing date values (raw serials, not display text).
$oldValues = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $oldValues += $ws.Cells.Item($r, $col).Value2
}

$count = $oldValues.Count

# Shift everything up by one row; the last row gets the next date
# (previous last value + 1) since there's no successor to pull from.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $i = $r - $firstRow
    if ($i -lt ($count - 1)) {
        $ws.Cells.Item($r, $col).Value = $oldValues[$i + 1]
    } else {
        $ws.Cells.Item($r, $col).Value = $oldValues[$count - 1] + 1
    }
}

# Match the selection left behind after performing the fill/shift over
# the whole date column.
$ws.Range("A2:A115").Select() | Out-Null
